$wb = $excel.ActiveWorkbook

$newPath = "C:\Katalon_mastercopy520\DataCommons_Automation\OutputFiles\TC02_Canine_Filter_Breed-AmerStaffd_Neo4jData.xlsx"

$wsMessage = $wb.Worksheets.Item("Message")
$wsMessage.Range("A10").Value = $newPath

$wsCypherMsg = $wb.Worksheets.Item("CypherOutput_Message")
$wsCypherMsg.Range("A10").Value = $newPath

$wsStatMsg = $wb.Worksheets.Item("StatOutput_Message")
$wsStatMsg.Range("A10").Value = $newPath
$wsStatMsg.Range("A20").Value = $newPath
